# Update spreadsheet tab labels.
#
# "Charter Schools" -> "Charter Schools in Ohio"
# "%Freq of School Type" -> "% Freq of School Type" (space added after %)
#
# Renaming a sheet via the .Name property makes Excel rewrite every
# formula elsewhere in the workbook that referenced the old tab name
# (the COUNTIF(...) formulas on "Frequencies" and "% Freq of School
# Type" both reference 'Charter Schools'!...), so those update for free.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Charter Schools").Name = "Charter Schools in Ohio"
$wb.Worksheets.Item("%Freq of School Type").Name = "% Freq of School Type"

# The remembered selection on "Charter Schools in Ohio" moved from G6 to
# F19. Selecting a range on a worksheet object (without Activate) doesn't
# disturb which tab/window is currently active.
$ws1 = $wb.Worksheets.Item("Charter Schools in Ohio")
$ws1.Range("F19").Select()

# The view on "% of Each Grade" was scrolled so row 7 becomes the top
# visible row (topLeftCell), while the selected cell remains D7.
$ws5 = $wb.Worksheets.Item("% of Each Grade")
$ws5.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws5.Range("D7").Select()

# "Charter Schools in Ohio" remains the workbook's selected/active tab
# (tabSelected="1"), so re-activate it last with its F19 selection intact.
$ws1.Activate()
$ws1.Range("F19").Select()
